$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab) from "Sheet1" to "Chetan Sakariya"
$ws.Name = "Chetan Sakariya"

# 2. Insert a new column before column A ("matchNo") — this shifts the
#    existing teamName..result columns from A:L to B:M
$ws.Columns.Item(1).Insert()

# 3. Insert a new row before row 2 — this shifts the existing data row
#    (originally row 2) down to row 3, making room for the new "51st" row
$ws.Rows.Item(2).Insert()

# Helper to fill a whole data row as TEXT (matches the workbook's existing
# convention of storing every value — including numeric-looking ones — as
# text, e.g. runs/balls/sr are plain text, not numbers).
function Set-TextRow($rowIndex, $values) {
    $rng = $ws.Range("A${rowIndex}:M${rowIndex}")
    $rng.NumberFormat = "@"
    $col = 1
    foreach ($v in $values) {
        $ws.Cells.Item($rowIndex, $col).Value = $v
        $col++
    }
}

# Header row (row 1)
Set-TextRow 1 @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

# Row 2 — new first data row
Set-TextRow 2 @("51st","Rajasthan Royals","Chetan Sakariya","b Coulter-Nile","6","11","1","0","54.54","Mumbai Indians","Sharjah","October 05","Mumbai won by 8 wickets (with 70 balls remaining)")

# Row 3 — previously existing row, now shifted down and given a matchNo
Set-TextRow 3 @("43rd","Rajasthan Royals","Chetan Sakariya","c de Villiers b Patel","2","2","0","0","100.00","Royal Challengers Bangalore","Dubai (DSC)","September 29","RCB won by 7 wickets (with 17 balls remaining)")

# Row 4 — new row
Set-TextRow 4 @("32nd","Rajasthan Royals","Chetan Sakariya","c & b Arshdeep Singh","7","6","1","0","116.66","Punjab Kings","Dubai (DSC)","September 21","Royals won by 2 runs")

# Row 5 — new row
Set-TextRow 5 @("16th","Rajasthan Royals","Chetan Sakariya","c †de Villiers b Patel","0","1","0","0","0.00","Royal Challengers Bangalore","Wankhede","April 22","RCB won by 10 wickets (with 21 balls remaining)")

# Row 6 — new row
Set-TextRow 6 @("54th","Rajasthan Royals","Chetan Sakariya","run out (Shakib Al Hasan/†Karthik)","1","5","0","0","20.00","Kolkata Knight Riders","Sharjah","October 07","KKR won by 86 runs")

# Row 7 — new row (states/D7 is empty)
Set-TextRow 7 @("12th","Rajasthan Royals","Chetan Sakariya","","0","0","0","0","-","Chennai Super Kings","Wankhede","April 19","Super Kings won by 45 runs")

Write-Output "edit applied"
